$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the stray _GoBack bookmark that sits on its own empty
# Heading1 paragraph (it gets re-created later, anchored to the new E23
# heading, in Change 3).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# Change 2: "3.What are the steps for creating product page in laravel 8
# ecommerce?" -> add the missing space after "3." and drop the (now
# unnecessary) proofErr gramStart/gramEnd wrapper around that run. Only the
# second "3.What" paragraph (the "product page" one) is touched; the first
# ("categories page") is left alone.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*creating product page in*") {
        $target = $p
        break
    }
}

$fixedXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DB353D" w:rsidRDefault="00DB353D" w:rsidP="00E140F0"><w:pPr><w:pStyle w:val="Heading1"/><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="030303"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="030303"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:rPr><w:t>3. What</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="030303"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:rPr><w:t xml:space="preserve"> are the steps for creating product page in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="030303"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:rPr><w:t>laravel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="030303"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:rPr><w:t xml:space="preserve"> 8 ecommerce?</w:t></w:r></w:p>
'@

$null = $target.Range.InsertXML($fixedXml)

# ---------------------------------------------------------------------------
# Change 3: the trailing blank paragraph at the very end of the document is
# replaced with the new "E23 - Admin show product categories in homepage?"
# section (heading + question + answer). The _GoBack bookmark now lives
# inside the heading run, right before "homepage?".
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$newSectionXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">E23 - </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Admin show product categories in </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>homepage?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>How to show Products Categories in homepage?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>First Create component and name it and then set the route for the given component.</w:t></w:r></w:p>
'@

$null = $lastPara.Range.InsertXML($newSectionXml)
